$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Losing Teams")

# Resize/extend the table (Table2) to include the two new columns first
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:AH97"))

# Header cells for the two new referee columns (AG = col 33, AH = col 34)
# Setting these after the resize causes the table's ListColumns to adopt
# the same names automatically (mirrors how Excel behaves interactively).
$ws.Cells.Item(1, 33).Value = "Referee 1"
$ws.Cells.Item(1, 34).Value = "Referee 2"

$refData = @{
    2 = @('Alessia Ferrari', 'Sebastien Dervieux')
    3 = @('Nora Debreceni', 'Frank Ohme')
    4 = @('Helene Painchaud', 'Liang Zhang')
    5 = @('Jennifer McCall', 'Aurely Blanchard')
    6 = @('Nora Debreceni', 'Nick Hodgers')
    7 = @('Adrian Alexandrescu', 'Andrej Franulovic')
    8 = @('Alessia Ferrari', 'Chisato Kurosaki')
    9 = @('Natalia Markopolou', 'Marta Cabanas')
    10 = @('Boris Margeta', 'Frank Ohme')
    11 = @('Helene Painchaud', 'Vojin Putnikovic')
    12 = @('Adrian Alexandrescu', 'Veselin Miskovic')
    13 = @('Alessia Ferrari', 'Georgios Stavridis')
    14 = @('Andrej Franulovic', 'Nora Debreceni')
    15 = @('Veselin Miskovic', 'Nick Hodgers')
    16 = @('Marta Cabanas', 'Aurely Blanchard')
    17 = @('Andrej Franulovic', 'Frank Ohme')
    18 = @('Jennifer McCall', 'Aurely Blanchard')
    19 = @('Vojin Putnikovic', 'Veselin Miskovic')
    20 = @('Alessia Ferrari', 'Sebastien Dervieux')
    21 = @('Vojin Putnikovic', 'Liang Zhang')
    22 = @('Aurely Blanchard', 'Natalia Markopolou')
    23 = @('Jennifer McCall', 'Chisato Kurosaki')
    24 = @('Marta Cabanas', 'Raffaele Colombo')
    25 = @('Nora Debreceni', 'Andrej Franulovic')
    26 = @('Jakov Blaskovic', 'Nicola Johnson')
    27 = @('Natalia Markopolou', 'Nikolett Sajben')
    28 = @('Matan Schwartz', 'Giuliana Nicolosi')
    29 = @('Danielle Dabbaghian', 'Frank Ohme')
    30 = @('Frank Ohme', 'Nicola Johnson')
    31 = @('Marta Cabanas', 'Jakov Blaskovic')
    32 = @('Matan Schwartz', 'Marieke van den Berg')
    33 = @('Danielle Dabbaghian', 'Nikolett Sajben')
    34 = @('Danielle Dabbaghian', 'Marta Cabanas')
    35 = @('Frank Ohme', 'Marieke van den Berg')
    36 = @('Jakov Blaskovic', 'Natalia Markopolou')
    37 = @('Marieke van den Berg', 'Marta Cabanas')
    38 = @('Nicola Johnson', 'Frank Ohme')
    39 = @('Jakov Blaskovic', 'Giuliana Nicolosi')
    40 = @('Matan Schwartz', 'Nikolett Sajben')
    41 = @('Marieke van den Berg', 'Natalia Markopolou')
    42 = @('Marta Cabanas', 'Jakov Blaskovic')
    43 = @('Matan Schwartz', 'Danielle Dabbaghian')
    44 = @('Giuliana Nicolosi', 'Frank Ohme')
    45 = @('Nicola Johnson', 'Marta Cabanas')
    46 = @('Natalia Markopolou', 'Jakov Blaskovic')
    47 = @('Matan Schwartz', 'Giuliana Nicolosi')
    48 = @('Frank Ohme', 'Nikolett Sajben')
    49 = @('Ruben Sap', 'Yuriko Udagawa')
    50 = @('Yang Peng', 'Ash Kaesler')
    51 = @('Julien Bourges', 'Nora Debreceni')
    52 = @('Jennifer McCall', 'Alessandro Severo')
    53 = @('Georgios Kravaritis', 'Ruben Sap')
    54 = @('Marta Cabanas', 'Nora Debreceni')
    55 = @('Alessandro Severo', 'Jennifer McCall')
    56 = @('Julien Bourges', 'Ash Kaesler')
    57 = @('Jennifer McCall', 'Ash Kaesler')
    58 = @('Georgios Kravaritis', 'Nora Debreceni')
    59 = @('Alessandro Severo', 'Marta Cabanas')
    66 = @('Fiona Haigh', 'Andrew Cairney')
    67 = @('Fiona Haigh', 'Andrew Cairney')
    68 = @('Marta Cabanas', 'Chisato Kurosaki')
    69 = @('Alessia Ferrari', 'Megan Rose Perry')
    70 = @('Fiona Haigh', 'Liang Zhang')
    71 = @('Aurely Blanchard', 'Tamas Kovacs Csatlos')
    72 = @('Marcella Braga', 'Boris Margeta')
    73 = @('Jennifer McCall', 'Frank Ohme')
    74 = @('German Moller', 'Megan Rose Perry')
    75 = @('Alessia Ferrari', 'Ivan Rakovic')
    76 = @('Jennifer McCall', 'Fiona Haigh')
    77 = @('David Gomez Pordomingo', 'Frank Ohme')
    78 = @('Aurely Blanchard', 'Boris Margeta')
    79 = @('Zhekang Wu', 'Scott Voltz')
    80 = @('Andrej Franulovic', 'Yasser Ali')
    81 = @('Tamas Kovacs Csatlos', 'Marta Cabanas')
    82 = @('Jennifer McCall', 'Marta Cabanas')
    83 = @('Zhekang Wu', 'German Moller')
    84 = @('Dasch Barber', 'Marcella Braga')
    85 = @('Aurely Blanchard', 'Fiona Haigh')
    86 = @('Liang Zhang', 'Michiel Zwart')
    87 = @('Alessia Ferrari', 'Marta Cabanas')
    88 = @('Jennifer McCall', 'Boris Margeta')
    89 = @('Aurely Blanchard', 'Maxim Gerasimov')
    90 = @('Jennifer McCall', 'German Moller')
    91 = @('Boris Margeta', 'Frank Ohme')
    92 = @('Tamas Kovacs Csatlos', 'Marta Cabanas')
    93 = @('Alessia Ferrari', 'Andrej Franulovic')
    94 = @('Aurely Blanchard', 'Nick Hodgers')
    95 = @('Andrej Franulovic', 'Michiel Zwart')
    96 = @('Ivan Rakovic', 'Jennifer McCall')
    97 = @('Alessia Ferrari', 'Marta Cabanas')
}

foreach ($rowNum in $refData.Keys) {
    $pair = $refData[$rowNum]
    $ws.Cells.Item([int]$rowNum, 33).Value = $pair[0]
    $ws.Cells.Item([int]$rowNum, 34).Value = $pair[1]
}

# Set the column widths for the new columns like the author did.
# (COM ColumnWidth values are offset by ~5/6 of a character versus the
# raw OOXML "width" units in this runtime, so compensate to land on 19.)
$ws.Columns.Item(33).ColumnWidth = 18.166666666666668
$ws.Columns.Item(34).ColumnWidth = 18.166666666666668

# Update selection to mirror the author's last active cell
$ws.Range("AG26").Select()
